$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "okay"
$ws.Range("E9").Value = "disabled"
